# implement worksheet_formula for xlsx
# Sheet1 (2nd tab, already active/selected) gets a new formula in A2 that
# adds B1 to the defined name OneRange, and the selection moves to A3.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Formula = "=B1+OneRange"
$ws.Range("A3").Select()
